# PC.xlsx update — "Add files via upload"
# Adds a new "GTX 1660" row to the "Placa de video" sheet, fills in two
# previously-zeroed "Preço Kabum" values on other sheets, and moves the
# active-sheet/selection state around (final active sheet becomes "Resumo").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Placa de video — fix I4 (Preço Kabum was 0, now matches Preço Atual)
#    and append a new product row (row 6: Asus GTX 1660).
# ---------------------------------------------------------------------
$wsPlaca = $wb.Worksheets.Item("Placa de video")

$wsPlaca.Range("I4").Value = 1550
$wsPlaca.Range("I4").NumberFormat = """R$""#,##0.00"

$wsPlaca.Range("B6").Value = 3

# Register the new URL as a shared string before "GTX 1660" so the
# workbook's shared-string table keeps the same ordering as upstream.
$gtx1660Url = "https://www.buscape.com.br/placa-de-video/placa-de-video-nvidia-geforce-gtx-1660-6-gb-gddr5-192-bits-asus-ph-gtx1660-o6g?_lc=88&q=gtx%201660"
[void]$wsPlaca.Hyperlinks.Add($wsPlaca.Cells.Item(6, 10), $gtx1660Url)
$wsPlaca.Range("J6").Style = "Hiperlink"

$wsPlaca.Range("C6").Value = "GTX 1660"
$wsPlaca.Range("D6").Value = "Asus"
$wsPlaca.Range("E6").Value = "6 GB GDDR5"
$wsPlaca.Range("F6").Value = "Kabum"

$wsPlaca.Range("G6").Value = 990
$wsPlaca.Range("H6").Value = 990
$wsPlaca.Range("I6").Value = 990
$wsPlaca.Range("G6:I6").NumberFormat = """R$""#,##0.00"

$wsPlaca.Activate()
[void]$wsPlaca.Range("I10").Select()

# ---------------------------------------------------------------------
# 2. Processador — only the remembered selection moved.
# ---------------------------------------------------------------------
$wsCpu = $wb.Worksheets.Item("Processador")
$wsCpu.Activate()
[void]$wsCpu.Range("D5").Select()

# ---------------------------------------------------------------------
# 3. Memória RAM — the 3rd kit's "Preço Kabum" (I6) was 0, now 630;
#    selection also moved.
# ---------------------------------------------------------------------
$wsRam = $wb.Worksheets.Item("Memória RAM")
$wsRam.Range("I6").Value = 630
$wsRam.Activate()
[void]$wsRam.Range("I9").Select()

# ---------------------------------------------------------------------
# 4. SSD — column F widened/best-fit; selection moved.
# ---------------------------------------------------------------------
$wsSsd = $wb.Worksheets.Item("SSD")
$wsSsd.Columns("F:F").ColumnWidth = 11
$wsSsd.Activate()
[void]$wsSsd.Range("H8").Select()

# ---------------------------------------------------------------------
# 5. Fonte — same column widening; selection moved.
# ---------------------------------------------------------------------
$wsFonte = $wb.Worksheets.Item("Fonte")
$wsFonte.Columns("F:F").ColumnWidth = 11
$wsFonte.Activate()
[void]$wsFonte.Range("D4").Select()

# ---------------------------------------------------------------------
# 6. Monitor — was the active tab; stays selected at F6 but loses the
#    "tabSelected" flag once another sheet is activated below.
# ---------------------------------------------------------------------
$wsMonitor = $wb.Worksheets.Item("Monitor")
$wsMonitor.Activate()
[void]$wsMonitor.Range("F6").Select()

# ---------------------------------------------------------------------
# 7. Resumo — becomes the new active tab (activeTab=8 / tabSelected).
# ---------------------------------------------------------------------
$wsResumo = $wb.Worksheets.Item("Resumo")
$wsResumo.Activate()
[void]$wsResumo.Range("H17").Select()
